$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:35:11.200472"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:11.200480"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:11.200483"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:11.200486"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:11.200489"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:11.200492"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:11.200494"
$dataSheet.Range("F9").Value = "2021-10-05 14:35:11.200497"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:11.200500"

# --- Add the new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row values
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row values
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Pancreatitis"
$metaSheet.Range("C2").Value = 154
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.3"
$metaSheet.Range("E2").Value = "2021-04-09T03:36:32.270823Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:11.196567"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/154/?format=json"

# Copy the header/index cell formatting from "data" sheet (same visual style: bold, bordered, centered)
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Range("A2").Value = 0

$dataSheet.Activate()
